$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New activity log entry: date + activity description
$ws.Range("A2").Value = 43777
$ws.Range("A2").NumberFormat = "d-mmm"
$ws.Range("B2").Value = "Découverte de Git et Github"

# Widen the description column so the text is readable (~74.57 chars, matches author's resize)
$ws.Columns.Item(2).ColumnWidth = 73.7

# Leave the selection where the author left it
$ws.Range("B11").Select() | Out-Null
